$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3:C4").EntireRow.Delete()
$ws.Range("A2").Value = "Gjuri bardhe"
$ws.Range("B2").Value = 41.451930028626599
$ws.Range("C2").Value = 20.0722199250377
$ws.Range("B2").Style = "Normal"
$ws.Rows(2).RowHeight = 30
$ws.Range("A2:XFD4").Select()
